$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add DNBC code values for "Mental health" (row 6) and "Loneliness" (row 8)
$ws.Range("F6").Value = "H052"
$ws.Range("F8").Value = '"H057", "H058", "H059"'

# Update the view: scroll position and selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F9").Select()
